# Mifos repayment schedule workbook - "modified test cases on overdue fix"
$wb = $excel.ActiveWorkbook

$wsInput   = $wb.Worksheets.Item(1)   # Input
$wsSummary = $wb.Worksheets.Item(2)   # Summary
$wsSched   = $wb.Worksheets.Item(3)   # Repayment schedule
$wsTrans   = $wb.Worksheets.Item(4)   # Transactions

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Sheet "Input"
# ---------------------------------------------------------------------------
$wsInput.Range("B2").Value = 42036

# ---------------------------------------------------------------------------
# Sheet "Summary" - restyle the percentage/amount cells to the plain style
# (style index 6 in the saved workbook) and refresh the overdue figures.
# ---------------------------------------------------------------------------
$plainRefSummary = $wsSummary.Range("A4")   # already carries the target style
$plainRefSummary.Copy()
$wsSummary.Range("B2").PasteSpecial($xlPasteFormats)
$wsSummary.Range("F2").PasteSpecial($xlPasteFormats)
$wsSummary.Range("A3:B3").PasteSpecial($xlPasteFormats)
$wsSummary.Range("E3:F3").PasteSpecial($xlPasteFormats)

$wsSummary.Range("G2").Clear()

$wsSummary.Range("A3").Value = 68.75
$wsSummary.Range("B3").Value = 12.5
$wsSummary.Range("E3").Value = 56.25
$wsSummary.Range("F3").Value = 13.89

# ---------------------------------------------------------------------------
# Sheet "Repayment schedule"
# ---------------------------------------------------------------------------
$plainRefSched = $wsSched.Range("A2")    # plain style (6)
$dateRefSched  = $wsSched.Range("C2")    # date style (7)
$italicRefSched = $wsTrans.Range("K2")   # italic/wrap style (10)

# Re-style the due/paid columns that drop their custom number formats
$plainRefSched.Copy()
$wsSched.Range("F3:F8").PasteSpecial($xlPasteFormats)
$wsSched.Range("H3:H8").PasteSpecial($xlPasteFormats)
$wsSched.Range("K3:K8").PasteSpecial($xlPasteFormats)
$wsSched.Range("L3").PasteSpecial($xlPasteFormats)
$wsSched.Range("M3").PasteSpecial($xlPasteFormats)
$wsSched.Range("P4:P8").PasteSpecial($xlPasteFormats)
$wsSched.Range("G7:G8").PasteSpecial($xlPasteFormats)

$dateRefSched.Copy()
$wsSched.Range("D3").PasteSpecial($xlPasteFormats)

$italicRefSched.Copy()
$wsSched.Range("E3").PasteSpecial($xlPasteFormats)

# Drop the now-empty "Over Due" helper column entries
$wsSched.Range("P2").Clear()
$wsSched.Range("O3:O8").Clear()

# Refresh the computed values row by row
$wsSched.Range("H3").Value = 12.5
$wsSched.Range("K3").Value = 845.83
$wsSched.Range("L3").Value = 845.83
$wsSched.Range("M3").Value = 0

$wsSched.Range("G4").Value = 3333.34
$wsSched.Range("H4").Value = 13.89
$wsSched.Range("K4").Value = 847.22
$wsSched.Range("P4").Value = 847.22

$wsSched.Range("G5").Value = 2500.0100000000002
$wsSched.Range("H5").Value = 17.36
$wsSched.Range("K5").Value = 850.69
$wsSched.Range("P5").Value = 850.69

$wsSched.Range("G6").Value = 1666.68

$wsSched.Range("G7").Value = 833.35
$wsSched.Range("H7").Value = 8.33
$wsSched.Range("K7").Value = 841.66
$wsSched.Range("P7").Value = 841.66

$wsSched.Range("G8").Value = 0.02

# New row 9 - clone formatting from the (now fully re-styled) row 8
$wsSched.Range("A8:P8").Copy()
$wsSched.Range("A9:P9").PasteSpecial($xlPasteFormats)

$wsSched.Range("A9").Value = 7
$wsSched.Range("B9").Value = 31
$wsSched.Range("C9").Value = 42217
$wsSched.Range("F9").Value = 0.02
$wsSched.Range("G9").Value = 0
$wsSched.Range("H9").Value = 0
$wsSched.Range("I9").Value = 0
$wsSched.Range("J9").Value = 0
$wsSched.Range("K9").Value = 0.02
$wsSched.Range("L9").Value = 0
$wsSched.Range("M9").Value = 0
$wsSched.Range("N9").Value = 0
$wsSched.Range("P9").Value = 0.02
$wsSched.Range("O9").Clear()

# ---------------------------------------------------------------------------
# Sheet "Transactions"
# ---------------------------------------------------------------------------
$plainRefTrans = $wsTrans.Range("A2")   # plain style (6)
$plainRefTrans.Copy()
$wsTrans.Range("F2").PasteSpecial($xlPasteFormats)
$wsTrans.Range("G2").PasteSpecial($xlPasteFormats)

$wsTrans.Range("A2").Value = 1705
$wsTrans.Range("E2").Value = 845.83
$wsTrans.Range("G2").Value = 12.5
$wsTrans.Range("J2").Value = 4166.67
$wsTrans.Range("A3").Value = 1513

# ---------------------------------------------------------------------------
# Selections - set last so the final active sheet/tab matches the target
# (Transactions ends up selected, matching activeTab=3 / tabSelected there).
# ---------------------------------------------------------------------------
$wsInput.Range("E10").Select()
$wsSummary.Range("E4").Select()
$wsSched.Range("J9").Select()
$wsTrans.Range("F3").Select()
